$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.466.67'
$ws.Range("E2").Value = '  +0.78%  '

$ws.Range("D3").Value = '1.724.73'
$ws.Range("E3").Value = '  +0.28%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.62'
$ws.Range("E5").Value = '  +2.05%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4822'
$ws.Range("E7").Value = '  +2.44%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2665'
$ws.Range("E8").Value = '  +1.62%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06178'
$ws.Range("E9").Value = '  -0.38%  '

$ws.Range("D10").Value = '1.736.36'
$ws.Range("E10").Value = '  +1.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07183'
$ws.Range("E11").Value = '  +1.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.55'
$ws.Range("E12").Value = '  +1.42%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6097'
$ws.Range("E13").Value = '  +2.09%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.521'
$ws.Range("E14").Value = '  +2.20%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '77.03'
$ws.Range("E15").Value = '  +1.11%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9999'
$ws.Range("E16").Value = '  -0.01%  '

$ws.Range("D17").Value = '26.479.24'
$ws.Range("E17").Value = '  +0.79%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.0000'
$ws.Range("E18").Value = '  -0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006944'
$ws.Range("E19").Value = '  +2.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.49'
$ws.Range("E20").Value = '  -0.60%  '

$ws.Range("D21").Value = '1.954.12'
$ws.Range("E21").Value = '  +0.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.515'
$ws.Range("E22").Value = '  -0.59%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.767'
$ws.Range("E23").Value = '  +0.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.238'
$ws.Range("E24").Value = '  -0.76%  '

$ws.Range("E25").Value = '  +1.66%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.30'
$ws.Range("E26").Value = '  +0.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.776'
$ws.Range("E27").Value = '  +0.88%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.394'
$ws.Range("E28").Value = '  -0.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '106.69'
$ws.Range("E29").Value = '  -0.68%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.960'
$ws.Range("E30").Value = '  -0.46%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08010'
$ws.Range("E31").Value = '  +3.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.682'
$ws.Range("E32").Value = '  +0.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04502'
$ws.Range("E33").Value = '  +1.10%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9997'
$ws.Range("E34").Value = '  -0.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.615'
$ws.Range("E35").Value = '  +0.03%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9924'
$ws.Range("E36").Value = '  +1.70%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6246'
$ws.Range("E37").Value = '  +1.10%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9123'
$ws.Range("E38").Value = '  -2.20%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.065'
$ws.Range("E39").Value = '  +7.21%  '

$ws.Range("E40").Value = '  -1.74%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.001'
$ws.Range("E41").Value = '  +0.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.33'
$ws.Range("E42").Value = '  -7.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.01502'
$ws.Range("E43").Value = '  +1.47%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.629'
$ws.Range("E44").Value = '  +3.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3858'
$ws.Range("E45").Value = '  +1.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.918'
$ws.Range("E46").Value = '  +10.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1178'
$ws.Range("E47").Value = '  -0.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05349'
$ws.Range("E48").Value = '  +1.66%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.35'
$ws.Range("E49").Value = '  +0.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.755'
$ws.Range("E50").Value = '  -0.36%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.250'
$ws.Range("E51").Value = '  +3.08%  '

